# Natmi following Dr Hou advice
# Update Dhh-Cdon LR-pair sheet: add "ECs" cluster and recompute rows for all
# Sending-cluster x Target-cluster combinations (ECs, FAPs, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Dhh", "Cdon", "ECs", 2, 0.6666666666666666, 2.179395, 6.538185, 0.4845018986408914, 0.4845018986408914, 2, 0.6666666666666666, 0.9894223333333333, 2.968267, 0.02008622060203138, 0.02008622060203139, 2.156342086155, 19.407078775395, 0.009731812018203995, 0.009731812018203997),
    @("ECs", "Dhh", "Cdon", "FAPs", 2, 0.6666666666666666, 2.179395, 6.538185, 0.4845018986408914, 0.4845018986408914, 3, 1, 34.79912266666667, 104.397368, 0.7064555054917406, 0.7064555054917406, 75.84103394412, 682.56930549708, 0.342279033716059, 0.342279033716059),
    @("ECs", "Dhh", "Cdon", "sCs", 2, 0.6666666666666666, 2.179395, 6.538185, 0.4845018986408914, 0.4845018986408914, 3, 1, 13.47021566666667, 40.410647, 0.2734582739062281, 0.2734582739062281, 29.356920672855, 264.212286055695, 0.1324910529066284, 0.1324910529066284),
    @("FAPs", "Dhh", "Cdon", "ECs", 3, 1, 0.7348883333333333, 2.204665, 0.1633732264179005, 0.1633732264179005, 2, 0.6666666666666666, 0.9894223333333333, 2.968267, 0.02008622060203138, 0.02008622060203139, 0.727114929506111, 6.544034365555, 0.003281550666295571, 0.003281550666295571),
    @("FAPs", "Dhh", "Cdon", "FAPs", 3, 1, 0.7348883333333333, 2.204665, 0.1633732264179005, 0.1633732264179005, 3, 1, 34.79912266666667, 104.397368, 0.7064555054917406, 0.7064555054917406, 25.57346925796889, 230.16122332172, 0.1154159152528745, 0.1154159152528745),
    @("FAPs", "Dhh", "Cdon", "sCs", 3, 1, 0.7348883333333333, 2.204665, 0.1633732264179005, 0.1633732264179005, 3, 1, 13.47021566666667, 40.410647, 0.2734582739062281, 0.2734582739062281, 9.899104340917221, 89.09193906825499, 0.04467576049873044, 0.04467576049873044),
    @("sCs", "Dhh", "Cdon", "ECs", 3, 1, 1.583934333333333, 4.751803, 0.3521248749412083, 0.3521248749412082, 2, 0.6666666666666666, 0.9894223333333333, 2.968267, 0.02008622060203138, 0.02008622060203139, 1.567180003933444, 14.104620035401, 0.007072857917531822, 0.007072857917531822),
    @("sCs", "Dhh", "Cdon", "FAPs", 3, 1, 1.583934333333333, 4.751803, 0.3521248749412083, 0.3521248749412082, 3, 1, 34.79912266666667, 104.397368, 0.7064555054917406, 0.7064555054917406, 55.11952516161156, 496.075726454504, 0.2487605565228072, 0.2487605565228072),
    @("sCs", "Dhh", "Cdon", "sCs", 3, 1, 1.583934333333333, 4.751803, 0.3521248749412083, 0.3521248749412082, 3, 1, 13.47021566666667, 40.410647, 0.2734582739062281, 0.2734582739062281, 21.33593707183789, 192.023433646541, 0.09629146050086923, 0.09629146050086922)
)

$startRow = 2
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Range("$($cols[$c])$r").Value = $rowVals[$c]
    }
}